# Edit script: applies row-data reordering (F:V columns) for matches that
# share an identical match date/time in the betexplorer Championship 2023-2024
# sheet, and appends the newly scraped Ipswich vs Millwall match as row 217.
#
# Columns A (Indice), B (pais), C (torneio), D (temporada) and E (data_partida)
# are positional/derived and are left untouched; only the match-specific data
# in columns F:V (home/away teams, scores, odds, timestamps, url) is swapped
# between rows, matching the permutation observed between the "before" and
# "after" canonical OOXML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cycle [125, 128]
$tmp = $ws.Range("F125:V125").Value()
$ws.Range("F125:V125").Value = $ws.Range("F128:V128").Value()
$ws.Range("F128:V128").Value = $tmp

# cycle [130, 132]
$tmp = $ws.Range("F130:V130").Value()
$ws.Range("F130:V130").Value = $ws.Range("F132:V132").Value()
$ws.Range("F132:V132").Value = $tmp

# cycle [131, 133]
$tmp = $ws.Range("F131:V131").Value()
$ws.Range("F131:V131").Value = $ws.Range("F133:V133").Value()
$ws.Range("F133:V133").Value = $tmp

# cycle [135, 143]
$tmp = $ws.Range("F135:V135").Value()
$ws.Range("F135:V135").Value = $ws.Range("F143:V143").Value()
$ws.Range("F143:V143").Value = $tmp

# cycle [136, 142]
$tmp = $ws.Range("F136:V136").Value()
$ws.Range("F136:V136").Value = $ws.Range("F142:V142").Value()
$ws.Range("F142:V142").Value = $tmp

# cycle [137, 141]
$tmp = $ws.Range("F137:V137").Value()
$ws.Range("F137:V137").Value = $ws.Range("F141:V141").Value()
$ws.Range("F141:V141").Value = $tmp

# cycle [138, 144, 140]
$tmp = $ws.Range("F138:V138").Value()
$ws.Range("F138:V138").Value = $ws.Range("F144:V144").Value()
$ws.Range("F144:V144").Value = $ws.Range("F140:V140").Value()
$ws.Range("F140:V140").Value = $tmp

# cycle [146, 147]
$tmp = $ws.Range("F146:V146").Value()
$ws.Range("F146:V146").Value = $ws.Range("F147:V147").Value()
$ws.Range("F147:V147").Value = $tmp

# cycle [157, 158]
$tmp = $ws.Range("F157:V157").Value()
$ws.Range("F157:V157").Value = $ws.Range("F158:V158").Value()
$ws.Range("F158:V158").Value = $tmp

# cycle [159, 164, 160, 163]
$tmp = $ws.Range("F159:V159").Value()
$ws.Range("F159:V159").Value = $ws.Range("F164:V164").Value()
$ws.Range("F164:V164").Value = $ws.Range("F160:V160").Value()
$ws.Range("F160:V160").Value = $ws.Range("F163:V163").Value()
$ws.Range("F163:V163").Value = $tmp

# cycle [161, 165, 162]
$tmp = $ws.Range("F161:V161").Value()
$ws.Range("F161:V161").Value = $ws.Range("F165:V165").Value()
$ws.Range("F165:V165").Value = $ws.Range("F162:V162").Value()
$ws.Range("F162:V162").Value = $tmp

# cycle [170, 171, 178, 174, 175, 179]
$tmp = $ws.Range("F170:V170").Value()
$ws.Range("F170:V170").Value = $ws.Range("F171:V171").Value()
$ws.Range("F171:V171").Value = $ws.Range("F178:V178").Value()
$ws.Range("F178:V178").Value = $ws.Range("F174:V174").Value()
$ws.Range("F174:V174").Value = $ws.Range("F175:V175").Value()
$ws.Range("F175:V175").Value = $ws.Range("F179:V179").Value()
$ws.Range("F179:V179").Value = $tmp

# cycle [172, 177]
$tmp = $ws.Range("F172:V172").Value()
$ws.Range("F172:V172").Value = $ws.Range("F177:V177").Value()
$ws.Range("F177:V177").Value = $tmp

# cycle [173, 176]
$tmp = $ws.Range("F173:V173").Value()
$ws.Range("F173:V173").Value = $ws.Range("F176:V176").Value()
$ws.Range("F176:V176").Value = $tmp

# cycle [184, 188]
$tmp = $ws.Range("F184:V184").Value()
$ws.Range("F184:V184").Value = $ws.Range("F188:V188").Value()
$ws.Range("F188:V188").Value = $tmp

# cycle [185, 186]
$tmp = $ws.Range("F185:V185").Value()
$ws.Range("F185:V185").Value = $ws.Range("F186:V186").Value()
$ws.Range("F186:V186").Value = $tmp

# cycle [197, 199]
$tmp = $ws.Range("F197:V197").Value()
$ws.Range("F197:V197").Value = $ws.Range("F199:V199").Value()
$ws.Range("F199:V199").Value = $tmp

# cycle [206, 209, 207, 210]
$tmp = $ws.Range("F206:V206").Value()
$ws.Range("F206:V206").Value = $ws.Range("F209:V209").Value()
$ws.Range("F209:V209").Value = $ws.Range("F207:V207").Value()
$ws.Range("F207:V207").Value = $ws.Range("F210:V210").Value()
$ws.Range("F210:V210").Value = $tmp

# cycle [213, 215]
$tmp = $ws.Range("F213:V213").Value()
$ws.Range("F213:V213").Value = $ws.Range("F215:V215").Value()
$ws.Range("F215:V215").Value = $tmp

# cycle [214, 216]
$tmp = $ws.Range("F214:V214").Value()
$ws.Range("F214:V214").Value = $ws.Range("F216:V216").Value()
$ws.Range("F216:V216").Value = $tmp

# --- append the new match row (217) ---

# Copy formatting (styles/number formats) from the last existing data row.
$ws.Range("A216:V216").Copy()
$ws.Range("A217:V217").PasteSpecial(-4122)

$ws.Range("A217").Value = 216
$ws.Range("B217").Value = "england"
$ws.Range("C217").Value = "championship"
$ws.Range("D217").Value = "2023-2024"
$ws.Range("E217").Value = 45259.875
$ws.Range("F217").Value = "Ipswich"
$ws.Range("G217").Value = 3
$ws.Range("H217").Value = "Millwall"
$ws.Range("I217").Value = 1
$ws.Range("J217").Value = 1.57
$ws.Range("K217").Value = "25/11/2023 18:42"
$ws.Range("L217").Value = 1.6
$ws.Range("M217").Value = "29/11/2023 20:40"
$ws.Range("N217").Value = 4.31
$ws.Range("O217").Value = "25/11/2023 18:42"
$ws.Range("P217").Value = 4.3
$ws.Range("Q217").Value = "29/11/2023 20:56"
$ws.Range("R217").Value = 5.66
$ws.Range("S217").Value = "25/11/2023 18:42"
$ws.Range("T217").Value = 5.82
$ws.Range("U217").Value = "29/11/2023 20:56"
$ws.Range("V217").Value = "https://www.betexplorer.com/football/england/championship/ipswich-millwall/8px1TzCA/"
